$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textLOT2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$textLOT2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"

# The "Requisitos" block lists two prerequisite courses in B23:C23 and
# B24:C24. The edit swaps their order: LOT2028 (Requisito fraco) now
# comes first (row 23), followed by LOT2038 (Indicação de Conjunto) in
# row 24.
$ws.Range("B23").Value = $textLOT2028
$ws.Range("C23").Value = $textLOT2028

$ws.Range("B24").Value = $textLOT2038
$ws.Range("C24").Value = $textLOT2038
